$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear stale content first (old Observação/Cobrança notes and the fully
# vacated rows 7-11) so the shared-string table gets compacted on save.
$ws.Range("E2:F6").ClearContents()
$ws.Range("H3:H6").ClearContents()
$ws.Range("A7:I11").ClearContents()

# Status column: every remaining ticket is now "Pendente".
$ws.Range("G2:G6").Value = "Pendente"

# Row 2 - MegaScan
$ws.Range("A2").Value = "Giovani"
$ws.Range("B2").Value = "'0643"
$ws.Range("C2").Value = "MegaScan"
$ws.Range("D2").Value = "Zona 5 aberta, cliente pedindo reparo."

# Row 3 - Viglioni
$ws.Range("A3").Value = "Giovani"
$ws.Range("B3").Value = "'0729"
$ws.Range("C3").Value = "Viglioni"
$ws.Range("D3").Value = "Sensor caiu no dia 22, no mesmo dia que foi recolocado. Pedi pra pregar de forma mais eficiente e realizar os devidos testes."

# Row 4 - Depósito Ideal
$ws.Range("A4").Value = "Giovani"
$ws.Range("B4").Value = "'0210"
$ws.Range("C4").Value = "Depósito Ideal"
$ws.Range("D4").Value = "Cliente solicitando revisão em um sensor e aumentar o tempo de saída."

# Row 5 - ViaMondo
$ws.Range("A5").Value = "Giovani"
$ws.Range("B5").Value = "'0774"
$ws.Range("C5").Value = "ViaMondo"
$ws.Range("D5").Value = "Pegar MAC da central e revisão nos setores abertos/disparo em falso."

# Row 6 - Rc Silva
$ws.Range("A6").Value = "Giovani"
$ws.Range("B6").Value = "'0355"
$ws.Range("C6").Value = "Rc Silva"
$ws.Range("D6").Value = "Zonas abertas, cliente pedindo reparo."

# Cobrança note, updated last
$ws.Range("H2").Value = "Maxvel: 38 / Forte: 17"

# Row heights: rows no longer holding long wrapped notes go back to the
# sheet's default (15pt, no explicit override); row 3 still wraps to two
# lines so it keeps its custom height.
$ws.Rows("2:11").EntireRow.AutoFit()
$ws.Rows(3).RowHeight = 30
